$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 51.59157666666666
$ws.Range("H2").Value = 154.77473
$ws.Range("I2").Value = 0.2641250550177587
$ws.Range("J2").Value = 0.2641250550177588
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 19.60726733333334
$ws.Range("N2").Value = 58.82180200000001
$ws.Range("O2").Value = 0.1509859438163708
$ws.Range("P2").Value = 0.1509859438163708
$ws.Range("Q2").Value = 1011.569835851496
$ws.Range("R2").Value = 9104.128522663459
$ws.Range("S2").Value = 0.03987917071740717
$ws.Range("T2").Value = 0.03987917071740718

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 51.59157666666666
$ws.Range("H3").Value = 154.77473
$ws.Range("I3").Value = 0.2641250550177587
$ws.Range("J3").Value = 0.2641250550177588
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 70.710031
$ws.Range("N3").Value = 212.130093
$ws.Range("O3").Value = 0.544503249041223
$ws.Range("P3").Value = 0.544503249041223
$ws.Range("Q3").Value = 3648.041985438877
$ws.Range("R3").Value = 32832.37786894989
$ws.Range("S3").Value = 0.1438169506103614
$ws.Range("T3").Value = 0.1438169506103615

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 51.59157666666666
$ws.Range("H4").Value = 154.77473
$ws.Range("I4").Value = 0.2641250550177587
$ws.Range("J4").Value = 0.2641250550177588
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 29.95517733333334
$ws.Range("N4").Value = 89.865532
$ws.Range("O4").Value = 0.2306701206736283
$ws.Range("P4").Value = 0.2306701206736284
$ws.Range("Q4").Value = 1545.434827956262
$ws.Range("R4").Value = 13908.91345160636
$ws.Range("S4").Value = 0.06092575831387513
$ws.Range("T4").Value = 0.06092575831387515

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.59157666666666
$ws.Range("H5").Value = 154.77473
$ws.Range("I5").Value = 0.2641250550177587
$ws.Range("J5").Value = 0.2641250550177588
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.589065333333332
$ws.Range("N5").Value = 28.767196
$ws.Range("O5").Value = 0.07384068646877778
$ws.Range("P5").Value = 0.0738406864687778
$ws.Range("Q5").Value = 494.7149993063421
$ws.Range("R5").Value = 4452.434993757079
$ws.Range("S5").Value = 0.019503175376115
$ws.Range("T5").Value = 0.01950317537611501

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09891538535728452
$ws.Range("J6").Value = 0.09891538535728453
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.60726733333334
$ws.Range("N6").Value = 58.82180200000001
$ws.Range("O6").Value = 0.1509859438163708
$ws.Range("P6").Value = 0.1509859438163708
$ws.Range("Q6").Value = 378.8350185949913
$ws.Range("R6").Value = 3409.515167354921
$ws.Range("S6").Value = 0.01493483281612963
$ws.Range("T6").Value = 0.01493483281612963

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09891538535728452
$ws.Range("J7").Value = 0.09891538535728453
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 70.710031
$ws.Range("N7").Value = 212.130093
$ws.Range("O7").Value = 0.544503249041223
$ws.Range("P7").Value = 0.544503249041223
$ws.Range("Q7").Value = 1366.199351155754
$ws.Range("R7").Value = 12295.79416040178
$ws.Range("S7").Value = 0.05385974870720604
$ws.Range("T7").Value = 0.05385974870720604

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09891538535728452
$ws.Range("J8").Value = 0.09891538535728453
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 29.95517733333334
$ws.Range("N8").Value = 89.865532
$ws.Range("O8").Value = 0.2306701206736283
$ws.Range("P8").Value = 0.2306701206736284
$ws.Range("Q8").Value = 578.7685743845246
$ws.Range("R8").Value = 5208.917169460721
$ws.Range("S8").Value = 0.02281682387684327
$ws.Range("T8").Value = 0.02281682387684328

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09891538535728452
$ws.Range("J9").Value = 0.09891538535728453
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.589065333333332
$ws.Range("N9").Value = 28.767196
$ws.Range("O9").Value = 0.07384068646877778
$ws.Range("P9").Value = 0.0738406864687778
$ws.Range("Q9").Value = 185.2718016286845
$ws.Range("R9").Value = 1667.44621465816
$ws.Range("S9").Value = 0.007303979957105579
$ws.Range("T9").Value = 0.007303979957105581

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 112.3724673333333
$ws.Range("H10").Value = 337.117402
$ws.Range("I10").Value = 0.5752951554216499
$ws.Range("J10").Value = 0.57529515542165
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.60726733333334
$ws.Range("N10").Value = 58.82180200000001
$ws.Range("O10").Value = 0.1509859438163708
$ws.Range("P10").Value = 0.1509859438163708
$ws.Range("Q10").Value = 2203.317007910934
$ws.Range("R10").Value = 19829.85307119841
$ws.Range("S10").Value = 0.08686148201432355
$ws.Range("T10").Value = 0.08686148201432356

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 112.3724673333333
$ws.Range("H11").Value = 337.117402
$ws.Range("I11").Value = 0.5752951554216499
$ws.Range("J11").Value = 0.57529515542165
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 70.710031
$ws.Range("N11").Value = 212.130093
$ws.Range("O11").Value = 0.544503249041223
$ws.Range("P11").Value = 0.544503249041223
$ws.Range("Q11").Value = 7945.860648686486
$ws.Range("R11").Value = 71512.74583817838
$ws.Range("S11").Value = 0.3132500812847637
$ws.Range("T11").Value = 0.3132500812847638

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 112.3724673333333
$ws.Range("H12").Value = 337.117402
$ws.Range("I12").Value = 0.5752951554216499
$ws.Range("J12").Value = 0.57529515542165
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 29.95517733333334
$ws.Range("N12").Value = 89.865532
$ws.Range("O12").Value = 0.2306701206736283
$ws.Range("P12").Value = 0.2306701206736284
$ws.Range("Q12").Value = 3366.137186354207
$ws.Range("R12").Value = 30295.23467718786
$ws.Range("S12").Value = 0.1327034029240657
$ws.Range("T12").Value = 0.1327034029240658

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 112.3724673333333
$ws.Range("H13").Value = 337.117402
$ws.Range("I13").Value = 0.5752951554216499
$ws.Range("J13").Value = 0.57529515542165
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.589065333333332
$ws.Range("N13").Value = 28.767196
$ws.Range("O13").Value = 0.07384068646877778
$ws.Range("P13").Value = 0.0738406864687778
$ws.Range("Q13").Value = 1077.546930927199
$ws.Range("R13").Value = 9697.92237834479
$ws.Range("S13").Value = 0.04248018919849683
$ws.Range("T13").Value = 0.04248018919849685

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 12.044915
$ws.Range("H14").Value = 36.134745
$ws.Range("I14").Value = 0.06166440420330686
$ws.Range("J14").Value = 0.06166440420330688
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 19.60726733333334
$ws.Range("N14").Value = 58.82180200000001
$ws.Range("O14").Value = 0.1509859438163708
$ws.Range("P14").Value = 0.1509859438163708
$ws.Range("Q14").Value = 236.1678684122767
$ws.Range("R14").Value = 2125.510815710491
$ws.Range("S14").Value = 0.00931045826851047
$ws.Range("T14").Value = 0.009310458268510473

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 12.044915
$ws.Range("H15").Value = 36.134745
$ws.Range("I15").Value = 0.06166440420330686
$ws.Range("J15").Value = 0.06166440420330688
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 70.710031
$ws.Range("N15").Value = 212.130093
$ws.Range("O15").Value = 0.544503249041223
$ws.Range("P15").Value = 0.544503249041223
$ws.Range("Q15").Value = 851.6963130423651
$ws.Range("R15").Value = 7665.266817381285
$ws.Range("S15").Value = 0.03357646843889184
$ws.Range("T15").Value = 0.03357646843889184

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 12.044915
$ws.Range("H16").Value = 36.134745
$ws.Range("I16").Value = 0.06166440420330686
$ws.Range("J16").Value = 0.06166440420330688
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 29.95517733333334
$ws.Range("N16").Value = 89.865532
$ws.Range("O16").Value = 0.2306701206736283
$ws.Range("P16").Value = 0.2306701206736284
$ws.Range("Q16").Value = 360.8075647899267
$ws.Range("R16").Value = 3247.26808310934
$ws.Range("S16").Value = 0.01422413555884419
$ws.Range("T16").Value = 0.01422413555884419

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 12.044915
$ws.Range("H17").Value = 36.134745
$ws.Range("I17").Value = 0.06166440420330686
$ws.Range("J17").Value = 0.06166440420330688
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.589065333333332
$ws.Range("N17").Value = 28.767196
$ws.Range("O17").Value = 0.07384068646877778
$ws.Range("P17").Value = 0.0738406864687778
$ws.Range("Q17").Value = 115.4994768694467
$ws.Range("R17").Value = 1039.49529182502
$ws.Range("S17").Value = 0.004553341937060365
$ws.Range("T17").Value = 0.004553341937060367
